$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Semester" column (E) for "DOB"
$ws.Columns("E:E").Insert()

# New column is text-formatted so the dd-mm-yyyy-style DOB strings are kept
# verbatim (not auto-converted to date serials).
$ws.Columns("E:E").NumberFormat = "@"

# Set the new column's header and data
$ws.Range("E1").Value = "DOB"
$ws.Range("E2").Value = "01-01-1998"
$ws.Range("E3").Value = "15-05-1999"
$ws.Range("E4").Value = "23-11-2000"

$ws.Columns("E:E").ColumnWidth = 12
